$d = $word.ActiveDocument

# The first (and only) paragraph currently reads "This is just for practice"
# followed by the _GoBack bookmark. We replace that paragraph's range with:
#   - the original sentence split into two runs ("This is just for practice" + ".")
#   - a brand-new second paragraph with the extra commentary, including the
#     proofErr spell-check markers around "git" and the moved _GoBack bookmark
#     at the very end (mirroring where Word leaves it after the last edit).
$p1 = $d.Paragraphs(1)
$target = $d.Range($p1.Range.Start, $p1.Range.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:r><w:t>This is just for practice</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:r><w:t xml:space="preserve">I had made comments to this doc but they were removed from the document by the command </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> restore &lt;filename&gt;</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$target.InsertXML($xml)
